# Contacts - Initial - 11 Oct 2024
# Adds new Contact-record fields (MiddleName, MailingStreet/City/State, Status,
# Office, Title, Department, LineOfBusiness) as columns H:P on the "Contact"
# sheet, for both existing sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

# ---- Header row (row 1) for the newly added columns H:P ----
$ws.Range("H1").Value = "MiddleName"
$ws.Range("I1").Value = "MailingStreet"
$ws.Range("J1").Value = "MailingCity"
$ws.Range("K1").Value = "MailingState"
$ws.Range("L1").Value = "Status"
$ws.Range("M1").Value = "Office"
$ws.Range("N1").Value = "Title"
$ws.Range("O1").Value = "Department"
$ws.Range("P1").Value = "LineOfBusiness"

# ---- Row 2 data (Sample John) ----
$ws.Range("H2").Value = "CK"
$ws.Range("I2").Value = "Street 1"
$ws.Range("J2").Value = "Kansas"
$ws.Range("K2").Value = "Missouri"
$ws.Range("L2").Value = "Active"
$ws.Range("M2").Value = "LA"
$ws.Range("N2").Value = "Associate"
$ws.Range("O2").Value = "CF"
$ws.Range("P2").Value = "CF"

# ---- Row 3 data (HRSample Jing) ----
$ws.Range("H3").Value = "CK"
$ws.Range("I3").Value = "Street 2"
$ws.Range("J3").Value = "Kansas"
$ws.Range("K3").Value = "Missouri"
$ws.Range("L3").Value = "Active"
$ws.Range("M3").Value = "LA"
$ws.Range("N3").Value = "Associate"
$ws.Range("O3").Value = "CF"
$ws.Range("P3").Value = "CF"

# Auto-fit the columns whose content is wider than the workbook default
# (mirrors the bestFit/customWidth columns Excel wrote for the wider fields).
$ws.Range("H1:K3").EntireColumn.AutoFit() | Out-Null
$ws.Range("O1:P3").EntireColumn.AutoFit() | Out-Null

# Reflect the new selection/scroll position recorded for the sheet view.
$ws.Range("O8").Select()
